$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Status column updates (was "Ready for handoff") ---
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# --- Latest Handback DateTime column (H) ---
$zhcn.Range("H2").Value = "2016-03-13 09:04:59"
$zhcn.Range("H3").Value = "2016-03-13 09:04:59"

$dede.Range("H2").Value = "2016-03-13 09:05:06"
$dede.Range("H3").Value = "2016-03-13 09:05:06"

Write-Output "status+datetime done"

# --- New "Latest Target File" (F) / "Latest Handback File" (G) columns ---
# zh-cn (sheet2) row 2 - 5fff7693 source file
$zhcn.Range("F2").Value = "5fff7693-9672-48a5-a098-040185dcd281.md"
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/940db9130cd23a455ec77d9b56bc0d191305caca/e2e/5fff7693-9672-48a5-a098-040185dcd281.md", "", "", "5fff7693-9672-48a5-a098-040185dcd281.md") | Out-Null

$zhcn.Range("G2").Value = "5fff7693-9672-48a5-a098-040185dcd281.49795b3c02d3501ff87a411f67f0c70fca8a892b.zh-cn.xlf"
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7bf4cd131e7678f00e2aef823110b5f50a7dfe95/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/5fff7693-9672-48a5-a098-040185dcd281.49795b3c02d3501ff87a411f67f0c70fca8a892b.zh-cn.xlf", "", "", "5fff7693-9672-48a5-a098-040185dcd281.49795b3c02d3501ff87a411f67f0c70fca8a892b.zh-cn.xlf") | Out-Null

# zh-cn (sheet2) row 3 - d38901d7 source file
$zhcn.Range("F3").Value = "d38901d7-938c-410b-b6cc-a1d01d19b6b0.md"
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/940db9130cd23a455ec77d9b56bc0d191305caca/e2e/d38901d7-938c-410b-b6cc-a1d01d19b6b0.md", "", "", "d38901d7-938c-410b-b6cc-a1d01d19b6b0.md") | Out-Null

$zhcn.Range("G3").Value = "d38901d7-938c-410b-b6cc-a1d01d19b6b0.0fa3a3fdfd0601929490cb78620888f992fa0d30.zh-cn.xlf"
$zhcn.Hyperlinks.Add($zhcn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7bf4cd131e7678f00e2aef823110b5f50a7dfe95/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d38901d7-938c-410b-b6cc-a1d01d19b6b0.0fa3a3fdfd0601929490cb78620888f992fa0d30.zh-cn.xlf", "", "", "d38901d7-938c-410b-b6cc-a1d01d19b6b0.0fa3a3fdfd0601929490cb78620888f992fa0d30.zh-cn.xlf") | Out-Null

# de-de (sheet3) row 2 - 5fff7693 source file
$dede.Range("F2").Value = "5fff7693-9672-48a5-a098-040185dcd281.md"
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/940db9130cd23a455ec77d9b56bc0d191305caca/e2e/5fff7693-9672-48a5-a098-040185dcd281.md", "", "", "5fff7693-9672-48a5-a098-040185dcd281.md") | Out-Null

$dede.Range("G2").Value = "5fff7693-9672-48a5-a098-040185dcd281.49795b3c02d3501ff87a411f67f0c70fca8a892b.de-de.xlf"
$dede.Hyperlinks.Add($dede.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d8592f88ee95f04f1657f9d3b259c77105383fdf/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/5fff7693-9672-48a5-a098-040185dcd281.49795b3c02d3501ff87a411f67f0c70fca8a892b.de-de.xlf", "", "", "5fff7693-9672-48a5-a098-040185dcd281.49795b3c02d3501ff87a411f67f0c70fca8a892b.de-de.xlf") | Out-Null

# de-de (sheet3) row 3 - d38901d7 source file
$dede.Range("F3").Value = "d38901d7-938c-410b-b6cc-a1d01d19b6b0.md"
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/940db9130cd23a455ec77d9b56bc0d191305caca/e2e/d38901d7-938c-410b-b6cc-a1d01d19b6b0.md", "", "", "d38901d7-938c-410b-b6cc-a1d01d19b6b0.md") | Out-Null

$dede.Range("G3").Value = "d38901d7-938c-410b-b6cc-a1d01d19b6b0.0fa3a3fdfd0601929490cb78620888f992fa0d30.de-de.xlf"
$dede.Hyperlinks.Add($dede.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d8592f88ee95f04f1657f9d3b259c77105383fdf/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d38901d7-938c-410b-b6cc-a1d01d19b6b0.0fa3a3fdfd0601929490cb78620888f992fa0d30.de-de.xlf", "", "", "d38901d7-938c-410b-b6cc-a1d01d19b6b0.0fa3a3fdfd0601929490cb78620888f992fa0d30.de-de.xlf") | Out-Null

Write-Output "hyperlink columns done"
